# Add the "email list" template content described by the commit:
#   - A1 and A2 both get the email address "xxx@email.com"
#   - Each cell becomes a mailto: hyperlink, which also applies Excel's
#     built-in "Hyperlink" cell style (underline + theme color 10 font)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$email = "xxx@email.com"

$ws.Range("A1").Value = $email
$ws.Range("A2").Value = $email

$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:$email")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:$email")
